$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 with the value that used to be in D4
$ws.Range("D2").Value = "Sell0 SpinaProcessunknown"

# Delete rows 3 through 6 (the other duplicate/obsolete rows)
$ws.Range("A3:D6").EntireRow.Delete()
